$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Row 10: Objetivos value (B/C) replaced by the "Docentes responsaveis"
# name that used to live at row 13. B10/C10 already exist, so the style
# (s=2 / s=3) is preserved automatically.
# ----------------------------------------------------------------------
$ws.Range("B10:C10").Value = "4894221 - Mariana Pereira de Melo"

# ----------------------------------------------------------------------
# Row 13: gains the "Programa resumido:" label in column A, and its
# B/C value becomes "Semestral".
# ----------------------------------------------------------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13:C13").Value = "Semestral"

# ----------------------------------------------------------------------
# Row 14: label becomes "Short syllabus:"; old long B/C text is removed
# (fully cleared so the cells no longer exist, matching the target).
# ----------------------------------------------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()

# ----------------------------------------------------------------------
# Row 15: label becomes "Programa:"; B/C gains "01/01/2021". Typing that
# literal text directly auto-converts it into a date serial, so instead
# enter it as a `="01/01/2021"` text formula, then collapse it down to a
# plain cached value via PasteSpecial values (still text, no formula
# left behind). B15/C15 didn't exist before, so finally paste formats
# from B8:C8 (same column styles, s=2 / s=3) to fix the style without
# registering any brand-new number format / style entries.
# ----------------------------------------------------------------------
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15:C15").Formula = "=""01/01/2021"""
$ws.Range("B15:C15").Copy() | Out-Null
$ws.Range("B15:C15").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = 0
$ws.Range("B8:C8").Copy() | Out-Null
$ws.Range("B15:C15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# Row 16: label becomes "Syllabus:"; old long B/C text is removed
# (fully cleared so the cells no longer exist, matching the target).
# ----------------------------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16:C16").Clear()

# ----------------------------------------------------------------------
# Row 17: label becomes "Avaliacao:".
# ----------------------------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"

# ----------------------------------------------------------------------
# Row 18: label becomes "Metodo:"; B/C gains the Mariana Pereira name.
# B18/C18 didn't exist before, so new cells default to the wrong style
# (column A's bold style); paste formats from B19:C19 (styles s=2/s=3)
# to fix that after setting the value.
# ----------------------------------------------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18:C18").Value = "4894221 - Mariana Pereira de Melo"
$ws.Range("B19:C19").Copy() | Out-Null
$ws.Range("B18:C18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# Row 19: label becomes "Criterio:" (B/C keep the NF=A... text as-is).
# ----------------------------------------------------------------------
$ws.Range("A19").Value = "Critério:"

# ----------------------------------------------------------------------
# Row 20: label becomes "Norma de recuperacao:" (B/C keep "NF>= 5,0.").
# ----------------------------------------------------------------------
$ws.Range("A20").Value = "Norma de recuperação:"

# ----------------------------------------------------------------------
# Row 21: label becomes "Bibliografia:" (B/C keep the "(NF+RC)/2..." text).
# ----------------------------------------------------------------------
$ws.Range("A21").Value = "Bibliografia:"

# ----------------------------------------------------------------------
# Row 22: label becomes "Requisitos:"; old bibliography B/C text removed
# (fully cleared so the cells no longer exist, matching the target).
# ----------------------------------------------------------------------
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").Clear()

# ----------------------------------------------------------------------
# Row 23: loses its "Requisitos:" label (moved to row 22, so A23 is fully
# cleared/removed) and gains the "LOB1012 - ..." requirement text in B/C
# (previously on row 24). B23/C23 didn't exist before, so paste formats
# from B21:C21 (styles s=2/s=3) after setting the value.
# ----------------------------------------------------------------------
$ws.Range("A23").Clear()
$ws.Range("B23:C23").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("B21:C21").Copy() | Out-Null
$ws.Range("B23:C23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# Row 24 now only duplicates row 23's old content; delete it entirely so
# the sheet shrinks back down to A1:C23.
# ----------------------------------------------------------------------
$ws.Rows.Item(24).Delete()

# ----------------------------------------------------------------------
# Row heights: rows 13,15,18 need 60/120/60; 21 needs 120; 23 needs 30;
# rows 17 and 22 drop back to the default (no custom height).
# ----------------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
